$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $value) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws "D2" "29.508.35"
Set-TextValue $ws "E2" "  -1.00%  "

Set-TextValue $ws "D3" "1.850.86"
Set-TextValue $ws "E3" "  -0.55%  "

Set-TextValue $ws "D4" "0.9983"
Set-TextValue $ws "E4" "  -0.08%  "

Set-TextValue $ws "D5" "241.92"
Set-TextValue $ws "E5" "  -0.78%  "

Set-TextValue $ws "D6" "0.6304"

Set-TextValue $ws "D7" "0.9994"
Set-TextValue $ws "E7" "  -0.02%  "

Set-TextValue $ws "D8" "0.07546"

Set-TextValue $ws "D9" "0.2980"
Set-TextValue $ws "E9" "  -0.32%  "

Set-TextValue $ws "D10" "24.43"
Set-TextValue $ws "E10" "  -1.18%  "

Set-TextValue $ws "D11" "0.07724"
Set-TextValue $ws "E11" "  +0.38%  "

Set-TextValue $ws "D12" "1.889.49"
Set-TextValue $ws "E12" "  +1.44%  "

Set-TextValue $ws "D13" "0.6929"
Set-TextValue $ws "E13" "  +0.05%  "

Set-TextValue $ws "D14" "5.016"
Set-TextValue $ws "E14" "  -0.77%  "

Set-TextValue $ws "D15" "83.72"
Set-TextValue $ws "E15" "  -0.48%  "

Set-TextValue $ws "D16" "0.000009781"
Set-TextValue $ws "E16" "  -1.25%  "

Set-TextValue $ws "D17" "2.156.39"
Set-TextValue $ws "E17" "  +1.63%  "

Set-TextValue $ws "D18" "6.246"
Set-TextValue $ws "E18" "  +2.20%  "

Set-TextValue $ws "D19" "29.552.32"
Set-TextValue $ws "E19" "  -0.86%  "

Set-TextValue $ws "D20" "234.46"
Set-TextValue $ws "E20" "  -0.82%  "

Set-TextValue $ws "D21" "12.51"
Set-TextValue $ws "E21" "  -1.20%  "

Set-TextValue $ws "D22" "0.9993"
Set-TextValue $ws "E22" "  -0.02%  "

Set-TextValue $ws "D23" "7.652"
Set-TextValue $ws "E23" "  +0.87%  "

Set-TextValue $ws "D24" "0.9992"
Set-TextValue $ws "E24" "  -0.10%  "

Set-TextValue $ws "D25" "154.80"
Set-TextValue $ws "E25" "  -2.31%  "

Set-TextValue $ws "D26" "0.1395"
Set-TextValue $ws "E26" "  -2.17%  "

Set-TextValue $ws "D27" "8.465"
Set-TextValue $ws "E27" "  -1.29%  "

Set-TextValue $ws "D28" "17.74"
Set-TextValue $ws "E28" "  -1.13%  "

Set-TextValue $ws "D29" "1.479"
Set-TextValue $ws "E29" "  -0.72%  "

Set-TextValue $ws "D30" "0.05876"
Set-TextValue $ws "E30" "  -5.01%  "

Set-TextValue $ws "D31" "1.251"
Set-TextValue $ws "E31" "  -2.86%  "

Set-TextValue $ws "D32" "4.110"
Set-TextValue $ws "E32" "  -1.21%  "

Set-TextValue $ws "D33" "4.054"
Set-TextValue $ws "E33" "  -1.06%  "

Set-TextValue $ws "D34" "1.881"
Set-TextValue $ws "E34" "  -0.55%  "

Set-TextValue $ws "D35" "1.171"
Set-TextValue $ws "E35" "  -0.30%  "

Set-TextValue $ws "D36" "0.7241"
Set-TextValue $ws "E36" "  -1.23%  "

Set-TextValue $ws "E37" "  -1.09%  "

Set-TextValue $ws "D38" "1.242.40"
Set-TextValue $ws "E38" "  +1.96%  "

Set-TextValue $ws "D39" "2.792"
Set-TextValue $ws "E39" "  -1.37%  "

Set-TextValue $ws "D40" "0.01785"
Set-TextValue $ws "E40" "  -0.46%  "

Set-TextValue $ws "D41" "0.9089"
Set-TextValue $ws "E41" "  -1.26%  "

Set-TextValue $ws "D42" "6.173"
Set-TextValue $ws "E42" "  -2.44%  "

Set-TextValue $ws "D43" "2.064.95"
Set-TextValue $ws "E43" "  +1.72%  "

Set-TextValue $ws "E44" "  -0.03%  "

Set-TextValue $ws "D45" "101.97"
Set-TextValue $ws "E45" "  -0.01%  "

Set-TextValue $ws "D46" "67.34"
Set-TextValue $ws "E46" "  +0.30%  "

Set-TextValue $ws "D47" "7.396"
Set-TextValue $ws "E47" "  +9.71%  "

Set-TextValue $ws "D48" "0.4048"

Set-TextValue $ws "D49" "9.148"
Set-TextValue $ws "E49" "  -0.28%  "

Set-TextValue $ws "B50" "BabyDogeCoin"
Set-TextValue $ws "C50" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue $ws "D50" "0.00000000117"
Set-TextValue $ws "E50" "  -1.16%  "

Set-TextValue $ws "B51" "RenderToken"
Set-TextValue $ws "C51" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws "D51" "1.704"
Set-TextValue $ws "E51" "  +2.08%  "
